# Weekly fruit/vegetable price update.
# Inserts a new daily price record (row 48) for "Murcott" mandarins from
# "Provincia de Limarí", pushing the existing rows 48-138 down by one
# (new last row becomes 139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 48, shifting rows 48:138
# down to 49:139 (and the used range / dimension grows to A1:T139).
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new record.
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = 'Vega Monumental Concepción'
$ws.Range("C48").Value = 'Bíobío'
$ws.Range("D48").Value = 44791
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 'Fruta'
$ws.Range("G48").Value = 100102
$ws.Range("H48").Value = 'Cítricos'
$ws.Range("I48").Value = 100102004
$ws.Range("J48").Value = 'Mandarina'
$ws.Range("K48").Value = 'Murcott'
$ws.Range("L48").Value = 'Primera'
$ws.Range("M48").Value = 100
$ws.Range("N48").Value = 7000
$ws.Range("O48").Value = 7500
$ws.Range("P48").Value = 7250
$ws.Range("Q48").Value = '$/bandeja 10 kilos'
$ws.Range("R48").Value = 'Provincia de Limarí'
$ws.Range("S48").Value = 725
$ws.Range("T48").Value = 10
